# Applies the "Automatic update of files" commit:
# adds Knärotsbuffertlänk / Klagomålslänk / Klagomålsmaillänk /
# Tillsynsbegäranslänk / Tillsynsbegäransmaillänk HYPERLINK formulas
# for a handful of rows in the Logging_HEBY block of the Avverkningsanmälningar sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

function Set-HyperlinkFormula {
    param(
        [string]$CellRef,
        [string]$Url
    )
    $formula = '=HYPERLINK("' + $Url + '")'
    $ws.Range($CellRef).Formula = $formula
}

# Row 2 - A 31572-2023 (gets the Knärotsbuffertlänk too, since it has Knärot listed)
Set-HyperlinkFormula "U2" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31572-2023.png"
Set-HyperlinkFormula "V2" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31572-2023.docx"
Set-HyperlinkFormula "W2" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31572-2023.docx"
Set-HyperlinkFormula "X2" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31572-2023.docx"
Set-HyperlinkFormula "Y2" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31572-2023.docx"

# Row 3 - A 32292-2023
Set-HyperlinkFormula "V3" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32292-2023.docx"
Set-HyperlinkFormula "W3" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32292-2023.docx"
Set-HyperlinkFormula "X3" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32292-2023.docx"
Set-HyperlinkFormula "Y3" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32292-2023.docx"

# Row 4 - A 32299-2023
Set-HyperlinkFormula "V4" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32299-2023.docx"
Set-HyperlinkFormula "W4" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32299-2023.docx"
Set-HyperlinkFormula "X4" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32299-2023.docx"
Set-HyperlinkFormula "Y4" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32299-2023.docx"

# Row 5 - A 32785-2023
Set-HyperlinkFormula "V5" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32785-2023.docx"
Set-HyperlinkFormula "W5" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32785-2023.docx"
Set-HyperlinkFormula "X5" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32785-2023.docx"
Set-HyperlinkFormula "Y5" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32785-2023.docx"

# Row 33 - A 31569-2023 (gets the Knärotsbuffertlänk too)
Set-HyperlinkFormula "U33" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31569-2023.png"
Set-HyperlinkFormula "V33" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31569-2023.docx"
Set-HyperlinkFormula "W33" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31569-2023.docx"
Set-HyperlinkFormula "X33" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31569-2023.docx"
Set-HyperlinkFormula "Y33" "https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31569-2023.docx"
